$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K15").Value = -0.3352267436446591
$ws.Range("J16").Value = 0.01855976243503714
$ws.Range("I17").Value = -0.1296176279974082
$ws.Range("H18").Value = -0.2870636170015632
$ws.Range("G19").Value = 0.2135958395245076
$ws.Range("F20").Value = -0.06676204101096155
$ws.Range("E21").Value = 0.1052128168340501
$ws.Range("D22").Value = -0.2006497229122814
$ws.Range("C23").Value = 0.4116802297750048
$ws.Range("B24").Value = -0.2766911554241067
